$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 4 changes
$ws.Range("G4").Value = 2.63
$ws.Range("I4").Value = 2.9
$ws.Range("J4").Value = 3.4
$ws.Range("L4").Value = 3.6
$ws.Range("X4").Value = 12
$ws.Range("Y4").Value = 11
$ws.Range("AI4").Value = 11
$ws.Range("AJ4").Value = 29
$ws.Range("AK4").Value = 26
$ws.Range("AN4").Value = 4.5
$ws.Range("AX4").Value = 17
$ws.Range("BA4").Value = 81

# Row 5 changes
$ws.Range("N5").Value = 8
$ws.Range("AF5").Value = 51
$ws.Range("AQ5").Value = 51

# Row 7 changes
$ws.Range("O7").Value = 1.33
$ws.Range("P7").Value = 3.25
